# feat: add 2022-Q4 data
#
# Before:
#   Sheet1 "总计"    : summary table, row2 = 2022-Q3 totals
#   Sheet2 "2022-Q3" : per-fund holdings for 2022-Q3
#
# After:
#   Sheet1 "总计"    : row2 = 2022-Q4 totals (new), row3 = 2022-Q3 totals (pushed down)
#   Sheet2 "2022-Q4" : per-fund holdings for 2022-Q4 (new data, same identity/rId as old sheet2)
#   Sheet3 "2022-Q3" : per-fund holdings for 2022-Q3 (the old sheet2 content, moved to a new sheet)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a value into a cell so it is stored as literal TEXT even when
# it "looks like" a number (e.g. "000788", "0.43"). A plain .Value assignment
# lets the engine infer a numeric type (exactly like typing into Excel), so
# we instead build the text with a formula and paste-special the computed
# value back in as a literal - this keeps the cell's style untouched (no
# NumberFormat/quotePrefix side effects) while forcing text storage.
# ---------------------------------------------------------------------------
function Set-TextValue {
    param($Sheet, $Cell, $Text)
    $stage = $Sheet.Range("ZZ1")
    $stage.Formula = "=" + '"' + $Text + '"'
    $stage.Copy()
    $Sheet.Range($Cell).PasteSpecial(-4163)
    $stage.Clear()
}

# ---------------------------------------------------------------------------
# 1. Sheet1 "总计": push the existing 2022-Q3 row down to row 3, then
#    overwrite row 2 with the new 2022-Q4 totals.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("A2:D2").Copy($ws1.Range("A3:D3"))
$ws1.Range("A3").Value = 1

$ws1.Range("B2").Value = "2022-Q4"
$ws1.Range("C2").Value = 1
$ws1.Range("D2").Value = 0.01

# ---------------------------------------------------------------------------
# 2. Duplicate the existing "2022-Q3" sheet (keeps all its fund rows/styles)
#    to a new sheet placed right after it; that copy becomes the new
#    "2022-Q3" sheet, while the original sheet is repurposed into "2022-Q4".
# ---------------------------------------------------------------------------
$wsOld = $wb.Worksheets.Item(2)
$wsOld.Copy($null, $wsOld)
$wsQ3 = $wb.Worksheets.Item(3)

$wsQ4 = $wsOld
$wsQ4.Name = "2022-Q4"
$wsQ3.Name = "2022-Q3"

# ---------------------------------------------------------------------------
# 3. Replace the fund-holding rows on the "2022-Q4" sheet with the new data
#    (headers in row 1 are unchanged, so only rows 2+ are cleared/rewritten).
# ---------------------------------------------------------------------------
$usedRows = $wsQ4.UsedRange.Rows.Count
if ($usedRows -ge 2) {
    $wsQ4.Range("A2:H" + $usedRows).Clear()
}

$wsQ4.Range("A2").Value = 0
Set-TextValue $wsQ4 "B2" "000788"
Set-TextValue $wsQ4 "C2" "前海开源中国成长灵活配置混合"
Set-TextValue $wsQ4 "D2" "0.43"
Set-TextValue $wsQ4 "E2" "85.88"
Set-TextValue $wsQ4 "F2" "1.70"
Set-TextValue $wsQ4 "G2" "0.0073"
$wsQ4.Range("H2").Value = 6
